# Update marksheet figures: correct marks count and corr/total marks ratio
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: number of correct answers (Right count) changed 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total marks obtained changed 57 -> 95
$ws.Range("B12").Value = 95

# "Total" row: Corr/total marks display text changed "54/84" -> "95/140"
$ws.Range("E12").Value = "95/140"
